# Publication 0.2.0 preparation:
#  - bump Version metadata value
#  - bump Date metadata value
#  - insert a new "Jurisdiction" / "iso:code:3166:FR" row into the Metadata table
#    (right after the "Contact" row), pushing all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above the current row 11 ("Description") so the table grows
# from 20 to 21 data rows, and fill it with the new Jurisdiction property.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Match the formatting of the other data rows (the insert otherwise leaves the
# new row with default/blank formatting instead of the shared data-row style).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Update Version (row 3) and Date (row 8) values in place.
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"
